$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rename the header row: "<name>_old" -> "<name>_FV2210" (cols A-J, 1-10)
#    and "<name>_new" -> "<name>_FV2304" (cols L-U, 12-21). Column K ("diff")
#    stays the same.
# ---------------------------------------------------------------------------
$baseNames = @("Segmentname","Segmentgruppe","Segment","Datenelement","Segment ID","Code","Qualifier","Beschreibung","Bedingungsausdruck","Bedingung")

for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $baseNames[$i] + "_FV2210"
}
for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $baseNames[$i] + "_FV2304"
}

# ---------------------------------------------------------------------------
# 2) Add a Table (ListObject) over A1:U72.
#    To avoid Excel auto-generating a header dxf / table style (because the
#    header row already carries custom formatting), we build the table on a
#    throw-away blank range first (where no special dxf gets produced), then
#    resize it onto the real A1:U72 range, and finally rename the columns by
#    touching the header cells again (which are already correctly named from
#    step 1, so re-assigning is a no-op that forces the table to re-read the
#    column names from the sheet).
# ---------------------------------------------------------------------------
$ws.Range("AA1").Value = "TmpCol1"
$ws.Range("AB1").Value = "TmpCol2"
$ws.Range("AA2").Value = "tmp"
$ws.Range("AB2").Value = "tmp"

$tmpRange = $ws.Range("AA1:AB2")
$tbl = $ws.ListObjects.Add(1, $tmpRange, $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

$ws.Range("AA1:AB2").Clear()

$tbl.Resize($ws.Range("A1:U72"))

# Re-touch header cells so the table column names are (re)read from the
# worksheet - they already hold the correct final text from step 1.
for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $baseNames[$i] + "_FV2210"
}
$ws.Cells.Item(1, 11).Value = "diff"
for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $baseNames[$i] + "_FV2304"
}

# ---------------------------------------------------------------------------
# 3) Freeze the header row (row 1).
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
